$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "Periodo Mora" column (E) for rows 16..53 previously listed periods
# 2005..2306 in ascending order (row16=2005 ... row53=2306).
# The database is rebuilt so the newest periods come first:
# row16=2306 ... row53=2005 (i.e. the list is reversed).
$periods = @(2306,2305,2304,2303,2302,2301,2212,2211,2210,2209,2208,2207,2206,2205,2204,2203,2202,2201,2112,2111,2110,2109,2108,2107,2106,2105,2104,2103,2102,2101,2012,2011,2010,2009,2008,2007,2006,2005)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value2 = [string]$periods[$i]
}

# The "Fecha" value (F) that used to sit with the 2306 period (row 53) now
# belongs with the 2306 period at row 16, and vice versa for the 2005 value.
$f16 = $ws.Cells.Item(16, 6).Value2
$f53 = $ws.Cells.Item(53, 6).Value2
$ws.Cells.Item(16, 6).Value2 = $f53
$ws.Cells.Item(53, 6).Value2 = $f16

$wb.Save()
